# Applies the "Updated cryptos list ... with GitHub Actions" price/volume refresh.
# D = Price column, E = Volume(1h) column. Both are plain text cells in the sheet
# (prices use "." as a thousands AND decimal separator in the source data, so they
# must stay text, not be reinterpreted as numbers) - use a leading apostrophe so the
# COM layer writes a literal string instead of coercing numeric-looking text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $ws.Cells.Item($row, $col).Formula = "'" + $text
}

Set-TextCell 2 4 "26.058.12"
Set-TextCell 2 5 "  +0.58%  "
Set-TextCell 3 4 "1.645.05"
Set-TextCell 3 5 "  +0.34%  "
Set-TextCell 4 4 "1.002"
Set-TextCell 4 5 "  -0.30%  "
Set-TextCell 5 4 "214.85"
Set-TextCell 5 5 "  -0.11%  "
Set-TextCell 6 4 "0.5099"
Set-TextCell 6 5 "  +1.31%  "
Set-TextCell 7 5 "  -0.32%  "
Set-TextCell 8 4 "0.2563"
Set-TextCell 8 5 "  -0.23%  "
Set-TextCell 9 4 "0.06358"
Set-TextCell 9 5 "  -0.54%  "
Set-TextCell 10 4 "19.57"
Set-TextCell 10 5 "  -0.13%  "
Set-TextCell 11 4 "0.07740"
Set-TextCell 11 5 "  -0.78%  "
Set-TextCell 12 4 "4.293"
Set-TextCell 12 5 "  +0.44%  "
Set-TextCell 13 4 "1.642.80"
Set-TextCell 13 5 "  -0.12%  "
Set-TextCell 14 4 "0.5441"
Set-TextCell 14 5 "  +0.20%  "
Set-TextCell 15 4 "64.30"
Set-TextCell 15 5 "  -0.58%  "
Set-TextCell 16 4 "0.0₅7721"
Set-TextCell 16 5 "  -1.77%  "
Set-TextCell 17 4 "26.063.11"
Set-TextCell 17 5 "  +0.49%  "
Set-TextCell 18 5 "  -0.41%  "
Set-TextCell 19 4 "198.84"
Set-TextCell 19 5 "  +0.24%  "
Set-TextCell 20 4 "4.433"
Set-TextCell 20 5 "  +1.12%  "
Set-TextCell 21 4 "9.927"
Set-TextCell 21 5 "  -0.15%  "
Set-TextCell 22 4 "6.051"
Set-TextCell 22 5 "  +1.31%  "
Set-TextCell 23 5 "  -0.33%  "
Set-TextCell 24 4 "1.871"
Set-TextCell 24 5 "  -0.32%  "
Set-TextCell 25 4 "141.13"
Set-TextCell 25 5 "  +0.84%  "
Set-TextCell 26 4 "0.1193"
Set-TextCell 26 5 "  +4.53%  "
Set-TextCell 27 5 "  -0.39%  "
Set-TextCell 28 4 "15.62"
Set-TextCell 28 5 "  -0.42%  "
Set-TextCell 29 5 "  -0.63%  "
Set-TextCell 30 4 "0.04859"
Set-TextCell 30 5 "  -0.34%  "
Set-TextCell 31 5 "  +0.03%  "
Set-TextCell 32 4 "3.167"
Set-TextCell 32 5 "  -0.73%  "
Set-TextCell 33 4 "1.528"
Set-TextCell 33 5 "  -0.17%  "
Set-TextCell 34 4 "2.366"
Set-TextCell 34 5 "  -0.10%  "
Set-TextCell 35 4 "0.8998"
Set-TextCell 35 5 "  +0.88%  "
Set-TextCell 36 4 "2.584"
Set-TextCell 36 5 "  -0.61%  "
Set-TextCell 37 4 "1.143.26"
Set-TextCell 37 5 "  +0.63%  "
Set-TextCell 38 4 "0.5470"
Set-TextCell 38 5 "  -1.33%  "
Set-TextCell 39 4 "0.01565"
Set-TextCell 39 5 "  +0.33%  "
Set-TextCell 40 5 "  -0.50%  "
Set-TextCell 41 4 "2.529"
Set-TextCell 41 5 "  -1.17%  "
Set-TextCell 42 4 "0.0₈130"
Set-TextCell 42 5 "  +4.55%  "
Set-TextCell 43 4 "0.8123"
Set-TextCell 43 5 "  -0.51%  "
Set-TextCell 44 4 "99.37"
Set-TextCell 44 5 "  -0.15%  "
Set-TextCell 45 4 "5.388"
Set-TextCell 45 5 "  -5.29%  "
Set-TextCell 46 4 "1.782.24"
Set-TextCell 46 5 "  +0.41%  "
Set-TextCell 47 4 "0.4525"
Set-TextCell 47 5 "  -0.05%  "
Set-TextCell 48 4 "54.96"
Set-TextCell 48 5 "  -0.62%  "
Set-TextCell 49 5 "  -0.73%  "
Set-TextCell 50 4 "0.05056"
Set-TextCell 50 5 "  -0.68%  "
Set-TextCell 51 4 "1.001"
Set-TextCell 51 5 "  -0.62%  "
